{"js": "// Remove the \"Appendix: Quick prototype\" sub-section: the Heading2\n// paragraph itself, the blank paragraph after it, the \"Figure: PDF\n// page 1\" caption paragraph, and the paragraph holding the embedded\n// page-1.png image (the block that sits between the \"Appendix: Links\"\n// summary paragraph and the \"Appendix: Links\" Heading2 paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].load(\"text,style\");\n}\nawait context.sync();\n\nlet startIndex = -1;\nlet endIndex = -1;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n\n  if (startIndex === -1 && text === \"Appendix: Quick prototype\" && para.style === \"Heading 2\") {\n    startIndex = i;\n  } else if (startIndex !== -1 && endIndex === -1 && text === \"Figure: PDF page 1\") {\n    // the very next paragraph holds the inline image and closes the block\n    endIndex = i + 1;\n  }\n}\n\nif (startIndex !== -1 && endIndex !== -1) {\n  for (let i = endIndex; i >= startIndex; i--) {\n    paragraphs.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the \"Appendix: Quick prototype\" sub-section: the Heading2\n# paragraph itself, the blank paragraph after it, the \"Figure: PDF\n# page 1\" caption paragraph, and the paragraph holding the embedded\n# page-1.png image \u2014 i.e. everything between the \"Appendix: Links\"\n# summary line and the \"Appendix: Links\" Heading2 that follows it.\n$d = $word.ActiveDocument\n\n$startIdx = -1\n$endIdx = -1\n$idx = 0\n\nforeach ($p in $d.Paragraphs) {\n    $idx = $idx + 1\n    $text = $p.Range.Text.Trim()\n\n    if ($startIdx -eq -1 -and $text -eq \"Appendix: Quick prototype\" -and $p.Style.NameLocal -eq \"Heading 2\") {\n        $startIdx = $idx\n    }\n    elseif ($startIdx -ne -1 -and $endIdx -eq -1 -and $text -eq \"Figure: PDF page 1\") {\n        # the very next paragraph holds the inline image and closes the block\n        $endIdx = $idx + 1\n    }\n}\n\nif ($startIdx -ne -1 -and $endIdx -ne -1) {\n    $startPar = $d.Paragraphs.Item($startIdx)\n    $endPar = $d.Paragraphs.Item($endIdx)\n    $rng = $d.Range($startPar.Range.Start, $endPar.Range.End)\n    $rng.Delete()\n}\n"}
